$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell E1 - same style (bold/centered) as the other header cells
$ws.Range("E1").Value = "Colocação"
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = -4108

# Ranking values for rows 2-8 (rows 9-10 stay unchanged / empty in column E)
$ws.Range("E2").Value = "1º"
$ws.Range("E3").Value = "2º"
$ws.Range("E4").Value = "3º"
$ws.Range("E5").Value = "4º"
$ws.Range("E6").Value = "5º"
$ws.Range("E7").Value = "6º"
$ws.Range("E8").Value = "21º"
